$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6472.647
$ws.Range("I62").Value = 2848.889
$ws.Range("K62").Value = 2848.889
$ws.Range("M62").Value = -2224.889
$ws.Range("H65").Value = 6472.647
$ws.Range("I65").Value = 2848.889
$ws.Range("K65").Value = 14244.445
$ws.Range("M65").Value = -11124.445
$ws.Range("H97").Value = 933.3333
$ws.Range("J97").Value = 900
$ws.Range("L97").Value = 2700
$ws.Range("N97").Value = -3692
$ws.Range("H99").Value = 3038.625
$ws.Range("I99").Value = 978.8
$ws.Range("J99").Value = 3974.9092
$ws.Range("K99").Value = 2936.4
$ws.Range("L99").Value = 11924.7276
$ws.Range("M99").Value = -1438.4
$ws.Range("N99").Value = -14920.7276
$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 14000
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 14000
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -15082
$ws.Range("H101").Value = 412.75
$ws.Range("I101").Value = 412.75
$ws.Range("K101").Value = 1238.25
$ws.Range("M101").Value = 383.75
$ws.Range("H103").Value = 125893.5
$ws.Range("J103").Value = 1021.2857
$ws.Range("L103").Value = 3063.8571
$ws.Range("N103").Value = -4235.8571
$ws.Range("H106").Value = 8337419.5
$ws.Range("I106").Value = 10004253
$ws.Range("K106").Value = 10004253
$ws.Range("M106").Value = -10003622
$ws.Range("H137").Value = 1489.4884
$ws.Range("I137").Value = 1498.1111
$ws.Range("J137").Value = 1483.28
$ws.Range("K137").Value = 4494.3333
$ws.Range("L137").Value = 4449.84
$ws.Range("M137").Value = -1944.3333
$ws.Range("N137").Value = -9549.84
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1316358.1
$ws.Range("I32").Value = 1529028.8
$ws.Range("J32").Value = 22611.666
$ws.Range("K32").Value = 1529028.8
$ws.Range("L32").Value = 22611.666
$ws.Range("M32").Value = -1528741.8
$ws.Range("N32").Value = -23185.666
$ws.Range("H61").Value = 3397.617
$ws.Range("I61").Value = 1978.8182
$ws.Range("J61").Value = 4646.16
$ws.Range("K61").Value = 1978.8182
$ws.Range("L61").Value = 4646.16
$ws.Range("M61").Value = -1766.8182
$ws.Range("N61").Value = -5070.16
$ws.Range("H97").Value = 250425
$ws.Range("I97").Value = 250425
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 250425
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -249929
$ws.Range("N97").ClearContents()
$ws.Range("H110").Value = 39764.305
$ws.Range("I110").Value = 69721.62
$ws.Range("K110").Value = 69721.62
$ws.Range("M110").Value = -67676.62
$ws.Range("H132").Value = 2486018.8
$ws.Range("I132").Value = 4607.2666
$ws.Range("J132").Value = 4812342
$ws.Range("K132").Value = 13821.7998
$ws.Range("L132").Value = 14437026
$ws.Range("M132").Value = -11291.7998
$ws.Range("N132").Value = -14442086
$ws.Range("H136").Value = 3397.617
$ws.Range("I136").Value = 1978.8182
$ws.Range("J136").Value = 4646.16
$ws.Range("K136").Value = 5936.4546
$ws.Range("L136").Value = 13938.48
$ws.Range("M136").Value = -3386.4546
$ws.Range("N136").Value = -19038.48
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1674.9318
$ws.Range("I86").Value = 1499.9697
$ws.Range("J86").Value = 2199.818
$ws.Range("K86").Value = 1499.9697
$ws.Range("L86").Value = 2199.818
$ws.Range("M86").Value = -376.9697000000001
$ws.Range("N86").Value = -4445.818
$ws.Range("H89").Value = 1674.9318
$ws.Range("I89").Value = 1499.9697
$ws.Range("J89").Value = 2199.818
$ws.Range("K89").Value = 7499.8485
$ws.Range("L89").Value = 10999.09
$ws.Range("M89").Value = -1883.8485
$ws.Range("N89").Value = -22231.09
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 6254.6665
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 8882
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 26646
$ws.Range("M19").Value = -2826
$ws.Range("N19").Value = -26994
$ws.Range("H93").Value = 4939
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 4939
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 14817
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -18561
$ws.Range("H94").Value = 6140.7144
$ws.Range("I94").Value = 350
$ws.Range("J94").Value = 7105.8335
$ws.Range("K94").Value = 1050
$ws.Range("L94").Value = 21317.5005
$ws.Range("M94").Value = -374
$ws.Range("N94").Value = -22669.5005
$ws.Range("H95").Value = 14000
$ws.Range("J95").Value = 14000
$ws.Range("L95").Value = 42000
$ws.Range("N95").Value = -46118
$ws.Range("H96").Value = 8000
$ws.Range("J96").Value = 8000
$ws.Range("L96").Value = 24000
$ws.Range("N96").Value = -28118
$ws.Range("H97").Value = 1502
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1502
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4506
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -5498
$ws.Range("H98").Value = 211388.78
$ws.Range("I98").Value = 299.66666
$ws.Range("J98").Value = 316933.34
$ws.Range("K98").Value = 898.9999799999999
$ws.Range("L98").Value = 950800.02
$ws.Range("M98").Value = 599.0000200000001
$ws.Range("N98").Value = -953796.02
$ws.Range("H99").Value = 5739
$ws.Range("I99").Value = 2650
$ws.Range("J99").Value = 8828
$ws.Range("K99").Value = 7950
$ws.Range("L99").Value = 26484
$ws.Range("M99").Value = -5704
$ws.Range("N99").Value = -30976
$ws.Range("H100").Value = 8462.4
$ws.Range("J100").Value = 8462.4
$ws.Range("L100").Value = 25387.2
$ws.Range("N100").Value = -27009.2
$ws.Range("H101").Value = 11129
$ws.Range("J101").Value = 11129
$ws.Range("L101").Value = 33387
$ws.Range("N101").Value = -38255
$ws.Range("H102").Value = 4532
$ws.Range("J102").Value = 4532
$ws.Range("L102").Value = 13596
$ws.Range("N102").Value = -18464
$ws.Range("H103").Value = 4264.846
$ws.Range("I103").Value = 235
$ws.Range("J103").Value = 6783.5
$ws.Range("K103").Value = 705
$ws.Range("L103").Value = 20350.5
$ws.Range("M103").Value = 174
$ws.Range("N103").Value = -22108.5
$ws.Range("H104").Value = 4022.5715
$ws.Range("J104").Value = 4022.5715
$ws.Range("L104").Value = 12067.7145
$ws.Range("N104").Value = -17309.7145
$ws.Range("H106").Value = 9354.076999999999
$ws.Range("J106").Value = 9354.076999999999
$ws.Range("L106").Value = 28062.231
$ws.Range("N106").Value = -29954.231
$ws.Range("H122").Value = 2159.8462
$ws.Range("I122").Value = 379.41934
$ws.Range("J122").Value = 3783.1765
$ws.Range("K122").Value = 3414.77406
$ws.Range("L122").Value = 34048.5885
$ws.Range("M122").Value = -964.7740599999997
$ws.Range("N122").Value = -38948.5885
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 33501.75
$ws.Range("J15").Value = 33501.75
$ws.Range("L15").Value = 33501.75
$ws.Range("N15").Value = -34077.75
$ws.Range("H54").Value = 8219
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 8219
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 8219
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -9259
$ws.Range("H81").Value = 5396.5557
$ws.Range("I81").Value = 5961.6665
$ws.Range("J81").Value = 4266.3335
$ws.Range("K81").Value = 11923.333
$ws.Range("L81").Value = 8532.666999999999
$ws.Range("M81").Value = -10862.333
$ws.Range("N81").Value = -10654.667
$ws.Range("H84").Value = 5396.5557
$ws.Range("I84").Value = 5961.6665
$ws.Range("J84").Value = 4266.3335
$ws.Range("K84").Value = 59616.665
$ws.Range("L84").Value = 42663.335
$ws.Range("M84").Value = -54312.665
$ws.Range("N84").Value = -53271.335
